# Automatische test-sync: 2025-06-19 21:23:50
# Append the new mail-log entry to "Logs" and refresh the "Dashboard"
# category summary (the new mail bumps "Offerte / Prijsaanvraag" from 1 to 2,
# which re-sorts it above the other single-count categories).

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 14 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A14").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D14").Value = "Offerte / Prijsaanvraag"
$logs.Range("F14").Value = "2025-06-19 21:23:17"
$logs.Range("G14").Value = "Nee"

# Extend the conditional formatting ranges from row 13 to row 14 so the new
# row picks up the same category / answered-status highlighting.
$catRules = $logs.Range("D2:D13").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D14"))
}

$answeredRules = $logs.Range("G2:G13").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G14"))
}

# --- Dashboard sheet: refresh category counts / ordering ------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Offerte / Prijsaanvraag"
$dash.Range("B4").Value = 2
$dash.Range("A6").Value = "Openingstijden / Locatie"
$dash.Range("A7").Value = "IT / Technisch probleem"
$dash.Range("A8").Value = "Sollicitatie / Vacature"
